$d = $word.ActiveDocument

# Locate the relevant paragraphs by their known text rather than trusting
# fixed indices: the "File: ..." line, the blank line right before it
# that precedes "Yield: ...", and the "Yield: ..." line itself.
$fileParaIdx = -1
$yieldParaIdx = -1
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("File: MOUNTA.90_X40")) {
        $fileParaIdx = $i
    }
    if ($t.StartsWith("Yield: 40")) {
        $yieldParaIdx = $i
    }
}
$blankParaIdx = $yieldParaIdx - 1

# --- Change 1: "File: MOUNTA.90_X40_..." -> "File: MOUNTA.90_X42_..." -------
# Word splits the run into three runs when the "0" in the middle of the
# existing run is retyped as "2" (the untouched text before/after the
# edited character stays in its own run). Reproduce that by: (a) fixing up
# the single changed character's text, then (b) nudging formatting on each
# of the three resulting segments (toggled back off) so Word keeps them as
# separate runs instead of silently re-merging identical-looking runs.
$paraRng = $d.Paragraphs.Item($fileParaIdx).Range
$paraStart = $paraRng.Start
$paraEnd = $paraRng.End

$fileRng = $d.Content
$fileRng.Find.Execute("X40_P2")
$digitStart = $fileRng.Start + 2
$digitEnd = $digitStart + 1

$digitRng = $d.Range($digitStart, $digitEnd)
$digitRng.Text = "2"

$seg1 = $d.Range($paraStart, $digitStart)
$seg1.Bold = 1
$seg1.Bold = 0

$seg2 = $d.Range($digitStart, $digitEnd)
$seg2.Bold = 1
$seg2.Bold = 0

$seg3 = $d.Range($digitEnd, $paraEnd)
$seg3.Bold = 1
$seg3.Bold = 0

# --- Change 2: empty run right before "Yield:" picks up the paragraph's
# East-Asian/complex-script font (Arial Unicode MS) ------------------------
$blankRng = $d.Paragraphs.Item($blankParaIdx).Range
$blankRng.Font.NameFarEast = "Arial Unicode MS"
$blankRng.Font.NameBi = "Arial Unicode MS"

# --- Change 3: "Yield: 40" -> "Yield: 42" ----------------------------------
$yieldRng = $d.Content
$yieldRng.Find.Execute("Yield: 40")
$yDigitEnd = $yieldRng.End
$yDigitStart = $yDigitEnd - 1

$yDigitRng = $d.Range($yDigitStart, $yDigitEnd)
$yDigitRng.Text = "2"

$yDigitRng2 = $d.Range($yDigitStart, $yDigitEnd)
$yDigitRng2.Bold = 1
$yDigitRng2.Bold = 0

Write-Output "done"
